$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45797.01041666666
$ws.Range("B2").Value = 766
$ws.Range("A3").Value = 45797.02083333334
$ws.Range("B3").Value = 733
$ws.Range("A4").Value = 45797.03125
$ws.Range("B4").Value = 762
$ws.Range("A5").Value = 45797.04166666666
$ws.Range("B5").Value = 809
$ws.Range("A6").Value = 45797.05208333334
$ws.Range("B6").Value = 895
$ws.Range("A7").Value = 45797.0625
$ws.Range("B7").Value = 950
$ws.Range("A8").Value = 45797.07291666666
$ws.Range("B8").Value = 1003
$ws.Range("A9").Value = 45797.08333333334
$ws.Range("B9").Value = 947
$ws.Range("A10").Value = 45797.09375
$ws.Range("B10").Value = 932
$ws.Range("A11").Value = 45797.10416666666
$ws.Range("B11").Value = 922
$ws.Range("A12").Value = 45797.11458333334
$ws.Range("B12").Value = 933
$ws.Range("A13").Value = 45797.125
$ws.Range("B13").Value = 1050
$ws.Range("A14").Value = 45797.13541666666
$ws.Range("B14").Value = 1078
$ws.Range("A15").Value = 45797.14583333334
$ws.Range("B15").Value = 987
$ws.Range("A16").Value = 45797.15625
$ws.Range("B16").Value = 936
$ws.Range("A17").Value = 45797.16666666666
$ws.Range("B17").Value = 903
$ws.Range("A18").Value = 45797.17708333334
$ws.Range("B18").Value = 829
$ws.Range("A19").Value = 45797.1875
$ws.Range("B19").Value = 833
$ws.Range("A20").Value = 45797.19791666666
$ws.Range("B20").Value = 830
$ws.Range("A21").Value = 45797.20833333334
$ws.Range("B21").Value = 804
$ws.Range("A22").Value = 45797.21875
$ws.Range("B22").Value = 784
$ws.Range("A23").Value = 45797.22916666666
$ws.Range("B23").Value = 698
$ws.Range("A24").Value = 45797.23958333334
$ws.Range("B24").Value = 600
$ws.Range("A25").Value = 45797.25
$ws.Range("B25").Value = 484
$ws.Range("A26").Value = 45797.26041666666
$ws.Range("B26").Value = 440
$ws.Range("A27").Value = 45797.27083333334
$ws.Range("B27").Value = 432
$ws.Range("A28").Value = 45797.28125
$ws.Range("B28").Value = 441
$ws.Range("A29").Value = 45797.29166666666
$ws.Range("B29").Value = 443
$ws.Range("A30").Value = 45797.30208333334
$ws.Range("B30").Value = 396
$ws.Range("A31").Value = 45797.3125
$ws.Range("B31").Value = 407
$ws.Range("A32").Value = 45797.32291666666
$ws.Range("B32").Value = 467
$ws.Range("A33").Value = 45797.33333333334
$ws.Range("B33").Value = 587
$ws.Range("A34").Value = 45797.34375
$ws.Range("B34").Value = 686
$ws.Range("A35").Value = 45797.35416666666
$ws.Range("B35").Value = 783
$ws.Range("A36").Value = 45797.36458333334
$ws.Range("B36").Value = 979
$ws.Range("A37").Value = 45797.375
$ws.Range("B37").Value = 1148
$ws.Range("A38").Value = 45797.38541666666
$ws.Range("B38").Value = 1359
$ws.Range("A39").Value = 45797.39583333334
$ws.Range("B39").Value = 1536
$ws.Range("A40").Value = 45797.40625
$ws.Range("B40").Value = 1593
$ws.Range("A41").Value = 45797.41666666666
$ws.Range("B41").Value = 0
$ws.Range("A42").Value = 45797.42708333334
$ws.Range("B42").Value = 0
$ws.Range("A43").Value = 45797.4375
$ws.Range("B43").Value = 0
$ws.Range("A44").Value = 45797.44791666666
$ws.Range("B44").Value = 0
$ws.Range("A45").Value = 45797.45833333334
$ws.Range("B45").Value = 0
$ws.Range("A46").Value = 45797.46875
$ws.Range("B46").Value = 0
$ws.Range("A47").Value = 45797.47916666666
$ws.Range("B47").Value = 0
$ws.Range("A48").Value = 45797.48958333334
$ws.Range("B48").Value = 0
$ws.Range("A49").Value = 45797.5
$ws.Range("B49").Value = 0
$ws.Range("A50").Value = 45797.51041666666
$ws.Range("B50").Value = 0
$ws.Range("A51").Value = 45797.52083333334
$ws.Range("B51").Value = 0
$ws.Range("A52").Value = 45797.53125
$ws.Range("B52").Value = 0
$ws.Range("A53").Value = 45797.54166666666
$ws.Range("B53").Value = 0
$ws.Range("A54").Value = 45797.55208333334
$ws.Range("B54").Value = 0
$ws.Range("A55").Value = 45797.5625
$ws.Range("B55").Value = 0
$ws.Range("A56").Value = 45797.57291666666
$ws.Range("B56").Value = 0
$ws.Range("A57").Value = 45797.58333333334
$ws.Range("B57").Value = 0
$ws.Range("A58").Value = 45797.59375
$ws.Range("B58").Value = 0
$ws.Range("A59").Value = 45797.60416666666
$ws.Range("B59").Value = 0
$ws.Range("A60").Value = 45797.61458333334
$ws.Range("B60").Value = 0
$ws.Range("A61").Value = 45797.625
$ws.Range("B61").Value = 0
$ws.Range("A62").Value = 45797.63541666666
$ws.Range("B62").Value = 0
$ws.Range("A63").Value = 45797.64583333334
$ws.Range("B63").Value = 0
$ws.Range("A64").Value = 45797.65625
$ws.Range("B64").Value = 0
$ws.Range("A65").Value = 45797.66666666666
$ws.Range("B65").Value = 0
$ws.Range("A66").Value = 45797.67708333334
$ws.Range("B66").Value = 0
$ws.Range("A67").Value = 45797.6875
$ws.Range("B67").Value = 0
$ws.Range("A68").Value = 45797.69791666666
$ws.Range("B68").Value = 0
$ws.Range("A69").Value = 45797.70833333334
$ws.Range("B69").Value = 0
$ws.Range("A70").Value = 45797.71875
$ws.Range("B70").Value = 0
$ws.Range("A71").Value = 45797.72916666666
$ws.Range("B71").Value = 0
$ws.Range("A72").Value = 45797.73958333334
$ws.Range("B72").Value = 0
$ws.Range("A73").Value = 45797.75
$ws.Range("B73").Value = 0
$ws.Range("A74").Value = 45797.76041666666
$ws.Range("B74").Value = 0
$ws.Range("A75").Value = 45797.77083333334
$ws.Range("B75").Value = 0
$ws.Range("A76").Value = 45797.78125
$ws.Range("B76").Value = 0
$ws.Range("A77").Value = 45797.79166666666
$ws.Range("B77").Value = 0
$ws.Range("A78").Value = 45797.80208333334
$ws.Range("B78").Value = 0
$ws.Range("A79").Value = 45797.8125
$ws.Range("B79").Value = 0
$ws.Range("A80").Value = 45797.82291666666
$ws.Range("B80").Value = 0
$ws.Range("A81").Value = 45797.83333333334
$ws.Range("B81").Value = 0
$ws.Range("A82").Value = 45797.84375
$ws.Range("B82").Value = 0
$ws.Range("A83").Value = 45797.85416666666
$ws.Range("B83").Value = 0
$ws.Range("A84").Value = 45797.86458333334
$ws.Range("B84").Value = 0
$ws.Range("A85").Value = 45797.875
$ws.Range("B85").Value = 0
$ws.Range("A86").Value = 45797.88541666666
$ws.Range("B86").Value = 0
$ws.Range("A87").Value = 45797.89583333334
$ws.Range("B87").Value = 0
$ws.Range("A88").Value = 45797.90625
$ws.Range("B88").Value = 0
$ws.Range("A89").Value = 45797.91666666666
$ws.Range("B89").Value = 0
$ws.Range("A90").Value = 45797.92708333334
$ws.Range("B90").Value = 0
$ws.Range("A91").Value = 45797.9375
$ws.Range("B91").Value = 0
$ws.Range("A92").Value = 45797.94791666666
$ws.Range("B92").Value = 0
$ws.Range("A93").Value = 45797.95833333334
$ws.Range("B93").Value = 0
$ws.Range("A94").Value = 45797.96875
$ws.Range("B94").Value = 0
$ws.Range("A95").Value = 45797.97916666666
$ws.Range("B95").Value = 0
$ws.Range("A96").Value = 45797.98958333334
$ws.Range("B96").Value = 0
$ws.Range("A97").Value = 45798
$ws.Range("B97").Value = 0